$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column A, shifting all existing columns right by one
$ws.Columns.Item(1).Insert()

# Copy the formatting (styles) from the new column B (the former column A)
# onto the new column A, so the new index column matches the header/body styles.
$ws.Range("B1:B33").Copy() | Out-Null
$ws.Range("A1:A33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Set the new column A's width to 23 (stored OOXML width), matching the
# ColumnWidth -> stored-width offset of 0.8333333333333333 used by this workbook.
$ws.Columns.Item(1).ColumnWidth = 23 - 0.8333333333333333

# Header label for the new index column
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

# Populate the new index column with the per-row index values
$indexValues = @(3110, 3249, 3250, 3251, 3252, 3253, 3254, 3255, 3256, 3257, 3258, 3259, 3260, 3261, 3262, 3263, 3264, 3265, 3266, 3267, 3268, 3269, 3270, 3271, 3272, 3273, 3274, 3275, 3276, 3277, 3278, 3279)

for ($i = 0; $i -lt $indexValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $indexValues[$i]
}
